{"js": "// The document ends with an empty numbered-list paragraph (level 0, numId 1).\n// This edit fills that paragraph with new text (\"Component with inline style\n// and inline template\") and appends a run of new list paragraphs after it,\n// finishing again with a new trailing empty paragraph at level 1 (replacing\n// the role the old empty paragraph used to play).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The last paragraph in the document is the blank trailing bullet.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Give it the new heading text (level 0 is already correct - inherited).\nlastParagraph.insertText(\"Component with inline style and inline template\", \"Replace\");\nawait context.sync();\n\n// New sub-bullets (level 1) under that heading.\nlet cursor = lastParagraph;\n\ncursor = cursor.insertParagraph(\"Inline style (how to generate)\", \"After\");\ncursor.listItem.level = 1;\nawait context.sync();\n\ncursor = cursor.insertParagraph(\"Inline template (how to generate)\", \"After\");\ncursor.listItem.level = 1;\nawait context.sync();\n\ncursor = cursor.insertParagraph(\"Both \", \"After\");\ncursor.listItem.level = 1;\nawait context.sync();\n\n// New top-level heading (level 0).\ncursor = cursor.insertParagraph(\"Module\", \"After\");\ncursor.listItem.level = 0;\nawait context.sync();\n\n// New sub-bullet (level 1) under \"Module\".\ncursor = cursor.insertParagraph(\"What\\u2019s a module? \\u2013 a complete feature. (e.g., User auth)\", \"After\");\ncursor.listItem.level = 1;\nawait context.sync();\n\n// New trailing empty sub-bullet (level 1), taking over the role of the\n// document's old trailing blank paragraph.\ncursor = cursor.insertParagraph(\"\", \"After\");\ncursor.listItem.level = 1;\nawait context.sync();\n", "ps1": "# The document ends with an empty numbered-list paragraph (level 1 / ilvl 0,\n# numId 1). This edit fills that paragraph with new text (\"Component with\n# inline style and inline template\") and appends a run of new list\n# paragraphs after it, finishing again with a new trailing empty paragraph\n# at level 2 / ilvl 1 (replacing the role the old empty paragraph used to\n# play).\n$d = $word.ActiveDocument\n\n# The last paragraph in the document is the blank trailing bullet - give it\n# the new heading text (its level, 1, i.e. ilvl 0, is already correct).\n$last = $d.Paragraphs.Last\n$last.Range.InsertBefore(\"Component with inline style and inline template\")\n\n# New sub-bullet (level 2 / ilvl 1) under that heading.\n$last.Range.InsertParagraphAfter()\n$p = $d.Paragraphs.Last\n$p.Range.InsertBefore(\"Inline style (how to generate)\")\n$p.Range.ListFormat.ListLevelNumber = 2\n\n# Another sub-bullet (level 2 / ilvl 1).\n$p.Range.InsertParagraphAfter()\n$p = $d.Paragraphs.Last\n$p.Range.InsertBefore(\"Inline template (how to generate)\")\n$p.Range.ListFormat.ListLevelNumber = 2\n\n# Another sub-bullet (level 2 / ilvl 1).\n$p.Range.InsertParagraphAfter()\n$p = $d.Paragraphs.Last\n$p.Range.InsertBefore(\"Both \")\n$p.Range.ListFormat.ListLevelNumber = 2\n\n# New top-level heading (level 1 / ilvl 0).\n$p.Range.InsertParagraphAfter()\n$p = $d.Paragraphs.Last\n$p.Range.InsertBefore(\"Module\")\n$p.Range.ListFormat.ListLevelNumber = 1\n\n# New sub-bullet (level 2 / ilvl 1) under \"Module\".\n$p.Range.InsertParagraphAfter()\n$p = $d.Paragraphs.Last\n$p.Range.InsertBefore(\"What\u2019s a module? \u2013 a complete feature. (e.g., User auth)\")\n$p.Range.ListFormat.ListLevelNumber = 2\n\n# New trailing empty sub-bullet (level 2 / ilvl 1), taking over the role of\n# the document's old trailing blank paragraph.\n$p.Range.InsertParagraphAfter()\n$p = $d.Paragraphs.Last\n$p.Range.ListFormat.ListLevelNumber = 2\n"}
